$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Set Up" paragraph: rework the wording around the event selection
#    and the screen that appears afterwards.
# ---------------------------------------------------------------------
$quoteOpen  = [char]0x201C
$quoteClose = [char]0x201D
$oldPhrase  = "select the event " + $quoteOpen + "User Testing" + $quoteClose
$d.Content.Find.Execute($oldPhrase, $true, $false, $false, $false, $false, $true, 1, $false, `
    "enter a random number for the event", 2) | Out-Null

$d.Content.Find.Execute("get up the age selection screen.", $true, $false, $false, $false, $false, $true, 1, $false, `
    "get up the drawing screen.", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Insert new questionnaire items + a "Notes From Tester" section
#    right after "How difficult was it to submit your drawing?"
# ---------------------------------------------------------------------
$anchorText = "How difficult was it to submit your drawing?"
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $anchorText) {
        $anchorPara = $p
        break
    }
}

$newParas = @(
    @{ Text = "How was the colour selection?"; Style = "Normal" },
    @{ Text = "What colours should be added?"; Style = "Normal" },
    @{ Text = "Do the icons make sense?"; Style = "Normal" },
    @{ Text = "Notes From Tester"; Style = "Heading 2" },
    @{ Text = "The focus for the tester is to see what questions are asked during the test. What kind of instructions the participant may ask, and make note of any struggles that the participant may have with the application"; Style = "Normal" }
)

$insertAfter = $anchorPara.Range
foreach ($item in $newParas) {
    $insertAfter.InsertParagraphAfter()
    $insertAfter = $insertAfter.Next(4, 1)
    $insertAfter.Text = $item.Text
    $currentStyle = $insertAfter.Paragraphs.Item(1).Style.NameLocal
    if ($currentStyle -ne $item.Style) {
        $insertAfter.Style = $item.Style
    }
}

# ---------------------------------------------------------------------
# 3) "Admin Application" heading becomes "Admin Application:"
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Admin Application") {
        $p.Range.InsertAfter(":") | Out-Null
        break
    }
}

# ---------------------------------------------------------------------
# 4) Add an extra empty paragraph before the final empty paragraph.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastStart = $lastPara.Range.Start
$collapsed = $d.Range($lastStart, $lastStart)
$collapsed.InsertParagraphBefore() | Out-Null
